$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Image1 (paragraph containing ${coQuanTrucThuoc}) - bump the wp:extent /
#    a:ext drawing sizes slightly (cap co quan update).
# ---------------------------------------------------------------------------
$p3 = $d.Paragraphs.Item(3)
$r3 = $p3.Range
$xml3 = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:wp="http://schemas.openxmlformats.org/drawingml/2006/wordprocessingDrawing" xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main" xmlns:wps="http://schemas.microsoft.com/office/word/2010/wordprocessingShape" xmlns:mc="http://schemas.openxmlformats.org/markup-compatibility/2006"><w:pPr><w:pStyle w:val="Normal"/><w:spacing w:lineRule="auto" w:line="240" w:before="0" w:after="0"/><w:jc w:val="center"/><w:rPr></w:rPr></w:pPr><w:r><mc:AlternateContent><mc:Choice Requires="wps"><w:drawing><wp:anchor behindDoc="0" distT="0" distB="0" distL="114300" distR="114300" simplePos="0" locked="0" layoutInCell="1" allowOverlap="1" relativeHeight="3"><wp:simplePos x="0" y="0"/><wp:positionH relativeFrom="margin"><wp:posOffset>736600</wp:posOffset></wp:positionH><wp:positionV relativeFrom="paragraph"><wp:posOffset>165100</wp:posOffset></wp:positionV><wp:extent cx="1048385" cy="19685"/><wp:effectExtent l="0" t="0" r="0" b="0"/><wp:wrapNone/><wp:docPr id="1" name="Image1"/><a:graphic xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main"><a:graphicData uri="http://schemas.microsoft.com/office/word/2010/wordprocessingShape"><wps:wsp><wps:cNvSpPr/><wps:spPr><a:xfrm><a:off x="0" y="0"/><a:ext cx="1047600" cy="19080"/></a:xfrm><a:custGeom><a:avLst/><a:gdLst/><a:ahLst/><a:rect l="l" t="t" r="r" b="b"/><a:pathLst><a:path w="21600" h="21600"><a:moveTo><a:pt x="0" y="0"/></a:moveTo><a:lnTo><a:pt x="21600" y="21600"/></a:lnTo></a:path></a:pathLst></a:custGeom><a:noFill/><a:ln w="9360"><a:solidFill><a:srgbClr val="000000"/></a:solidFill><a:miter/></a:ln></wps:spPr><wps:style><a:lnRef idx="0"/><a:fillRef idx="0"/><a:effectRef idx="0"/><a:fontRef idx="minor"/></wps:style><wps:bodyPr/></wps:wsp></a:graphicData></a:graphic></wp:anchor></w:drawing></mc:Choice><mc:Fallback><w:pict/></mc:Fallback></mc:AlternateContent></w:r><w:r><w:rPr><w:rFonts w:eastAsia="Times New Roman" w:cs="Times New Roman" w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:b/><w:position w:val="0"/><w:sz w:val="26"/><w:sz w:val="26"/><w:szCs w:val="26"/><w:vertAlign w:val="baseline"/></w:rPr><w:t>${coQuanTrucThuoc}</w:t></w:r></w:p>'
$r3.InsertXML($xml3)

# ---------------------------------------------------------------------------
# 2) Image2 (paragraph containing "Độc lập – Tự do – Hạnh phúc") - same kind
#    of extent bump.
# ---------------------------------------------------------------------------
$p6 = $d.Paragraphs.Item(6)
$r6 = $p6.Range
$xml6 = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:wp="http://schemas.openxmlformats.org/drawingml/2006/wordprocessingDrawing" xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main" xmlns:wps="http://schemas.microsoft.com/office/word/2010/wordprocessingShape" xmlns:mc="http://schemas.openxmlformats.org/markup-compatibility/2006"><w:pPr><w:pStyle w:val="Normal"/><w:spacing w:lineRule="auto" w:line="240" w:before="0" w:after="0"/><w:jc w:val="center"/><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:eastAsia="Times New Roman" w:cs="Times New Roman"/><w:b w:val="false"/><w:b w:val="false"/><w:position w:val="0"/><w:sz w:val="22"/><w:sz w:val="26"/><w:szCs w:val="26"/><w:vertAlign w:val="baseline"/></w:rPr></w:pPr><w:r><mc:AlternateContent><mc:Choice Requires="wps"><w:drawing><wp:anchor behindDoc="0" distT="0" distB="0" distL="114300" distR="114300" simplePos="0" locked="0" layoutInCell="1" allowOverlap="1" relativeHeight="2"><wp:simplePos x="0" y="0"/><wp:positionH relativeFrom="margin"><wp:posOffset>143510</wp:posOffset></wp:positionH><wp:positionV relativeFrom="paragraph"><wp:posOffset>165100</wp:posOffset></wp:positionV><wp:extent cx="2051685" cy="19685"/><wp:effectExtent l="0" t="0" r="0" b="0"/><wp:wrapNone/><wp:docPr id="2" name="Image2"/><a:graphic xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main"><a:graphicData uri="http://schemas.microsoft.com/office/word/2010/wordprocessingShape"><wps:wsp><wps:cNvSpPr/><wps:spPr><a:xfrm><a:off x="0" y="0"/><a:ext cx="2050920" cy="19080"/></a:xfrm><a:custGeom><a:avLst/><a:gdLst/><a:ahLst/><a:rect l="l" t="t" r="r" b="b"/><a:pathLst><a:path w="21600" h="21600"><a:moveTo><a:pt x="0" y="0"/></a:moveTo><a:lnTo><a:pt x="21600" y="21600"/></a:lnTo></a:path></a:pathLst></a:custGeom><a:noFill/><a:ln w="9360"><a:solidFill><a:srgbClr val="000000"/></a:solidFill><a:miter/></a:ln></wps:spPr><wps:style><a:lnRef idx="0"/><a:fillRef idx="0"/><a:effectRef idx="0"/><a:fontRef idx="minor"/></wps:style><wps:bodyPr/></wps:wsp></a:graphicData></a:graphic></wp:anchor></w:drawing></mc:Choice><mc:Fallback><w:pict/></mc:Fallback></mc:AlternateContent></w:r><w:r><w:rPr><w:rFonts w:eastAsia="Times New Roman" w:cs="Times New Roman" w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:b/><w:position w:val="0"/><w:sz w:val="26"/><w:sz w:val="26"/><w:szCs w:val="26"/><w:vertAlign w:val="baseline"/></w:rPr><w:t>Độc lập – Tự do – Hạnh phúc</w:t></w:r></w:p>'
$r6.InsertXML($xml6)

# ---------------------------------------------------------------------------
# 3) ${capHanhChinh} gains 63 trailing spaces (right-padding placeholder so
#    the printed form lines up under "cap co quan").
# ---------------------------------------------------------------------------
$rFind = $d.Content
$rFind.Find.ClearFormatting()
$foundCap = $rFind.Find.Execute('${capHanhChinh}')
if ($foundCap) {
    $rFind.Collapse(0)
    $rFind.InsertAfter("                                                               ")
}

# ---------------------------------------------------------------------------
# 4) Merge the trailing " " + "2018" runs in the "Đà Nẵng, ngày ... năm 2018"
#    line into a single run (same formatting on both sides already).
# ---------------------------------------------------------------------------
$rFind2 = $d.Content
$rFind2.Find.ClearFormatting()
$rFind2.Find.Execute(" 2018", $true, $false, $false, $false, $false, $true, 1, $false, " 2018", 2)
